$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1369.5
$ws.Range("I28").Value = 453.1905
$ws.Range("K28").Value = 453.1905
$ws.Range("M28").Value = 31.80950000000001
$ws.Range("H76").Value = 11629.375
$ws.Range("I76").Value = 13279.091
$ws.Range("K76").Value = 13279.091
$ws.Range("M76").Value = -12964.091
$ws.Range("H79").Value = 11629.375
$ws.Range("I79").Value = 13279.091
$ws.Range("K79").Value = 13279.091
$ws.Range("M79").Value = -12187.091
$ws.Range("H80").Value = 569.5217
$ws.Range("I80").Value = 488
$ws.Range("J80").Value = 675.5
$ws.Range("K80").Value = 1464
$ws.Range("L80").Value = 2026.5
$ws.Range("M80").Value = -466
$ws.Range("N80").Value = -4022.5
$ws.Range("H83").Value = 569.5217
$ws.Range("I83").Value = 488
$ws.Range("J83").Value = 675.5
$ws.Range("K83").Value = 4392
$ws.Range("L83").Value = 6079.5
$ws.Range("M83").Value = 600
$ws.Range("N83").Value = -16063.5
$ws.Range("H137").Value = 2973.6956
$ws.Range("I137").Value = 3420.9375
$ws.Range("J137").Value = 1951.4286
$ws.Range("K137").Value = 10262.8125
$ws.Range("L137").Value = 5854.2858
$ws.Range("M137").Value = -7712.8125
$ws.Range("N137").Value = -10954.2858

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11369209
$ws.Range("I32").Value = 14287845
$ws.Range("K32").Value = 14287845
$ws.Range("M32").Value = -14287558
$ws.Range("H45").Value = 1755.875
$ws.Range("J45").Value = 2483.3333
$ws.Range("L45").Value = 2483.3333
$ws.Range("N45").Value = -3237.3333
$ws.Range("H63").Value = 4205.294
$ws.Range("I63").Value = 2711.7693
$ws.Range("K63").Value = 2711.7693
$ws.Range("M63").Value = -2025.7693
$ws.Range("H66").Value = 4205.294
$ws.Range("I66").Value = 2711.7693
$ws.Range("K66").Value = 13558.8465
$ws.Range("M66").Value = -10126.8465
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3900.9285
$ws.Range("I99").Value = 2322.6667
$ws.Range("K99").Value = 2322.6667
$ws.Range("M99").Value = -824.6667000000002
$ws.Range("H105").Value = 8354.450000000001
$ws.Range("I105").Value = 11726.818
$ws.Range("J105").Value = 4232.6665
$ws.Range("K105").Value = 11726.818
$ws.Range("L105").Value = 4232.6665
$ws.Range("M105").Value = -9979.817999999999
$ws.Range("N105").Value = -7726.6665
$ws.Range("H107").Value = 3018.55
$ws.Range("I107").Value = 2492.4707
$ws.Range("J107").Value = 5999.6665
$ws.Range("K107").Value = 2492.4707
$ws.Range("L107").Value = 5999.6665
$ws.Range("M107").Value = -572.4706999999999
$ws.Range("N107").Value = -9839.666499999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H62").Value = 4840.3335
$ws.Range("I62").Value = 3839.8
$ws.Range("J62").Value = 5555
$ws.Range("K62").Value = 3839.8
$ws.Range("L62").Value = 5555
$ws.Range("M62").Value = -3215.8
$ws.Range("N62").Value = -6803
$ws.Range("H65").Value = 4840.3335
$ws.Range("I65").Value = 3839.8
$ws.Range("J65").Value = 5555
$ws.Range("K65").Value = 19199
$ws.Range("L65").Value = 27775
$ws.Range("M65").Value = -16079
$ws.Range("N65").Value = -34015
$ws.Range("H140").Value = 101700
$ws.Range("J140").Value = 101700
$ws.Range("L140").Value = 101700
$ws.Range("N140").Value = -112060

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 158567.42
$ws.Range("J37").Value = 158567.42
$ws.Range("L37").Value = 475702.26
$ws.Range("N37").Value = -475926.26
$ws.Range("H39").Value = 3184.6365
$ws.Range("J39").Value = 3999.875
$ws.Range("L39").Value = 11999.625
$ws.Range("N39").Value = -12587.625
$ws.Range("H56").Value = 21211
$ws.Range("I56").Value = 21211
$ws.Range("K56").Value = 21211
$ws.Range("M56").Value = -20681
$ws.Range("H88").Value = 10999.571
$ws.Range("J88").Value = 10999.571
$ws.Range("L88").Value = 32998.713
$ws.Range("N88").Value = -33854.713
$ws.Range("H91").Value = 10999.571
$ws.Range("J91").Value = 10999.571
$ws.Range("L91").Value = 32998.713
$ws.Range("N91").Value = -35962.713
$ws.Range("H113").Value = 3263.3333
$ws.Range("J113").Value = 3384.5557
$ws.Range("L113").Value = 10153.6671
$ws.Range("N113").Value = -14493.6671
$ws.Range("H122").Value = 1697.8334
$ws.Range("I122").Value = 1549.5
$ws.Range("J122").Value = 1772
$ws.Range("K122").Value = 13945.5
$ws.Range("L122").Value = 15948
$ws.Range("M122").Value = -11495.5
$ws.Range("N122").Value = -20848
$ws.Range("H129").Value = 4081.8667
$ws.Range("J129").Value = 3948.0454
$ws.Range("L129").Value = 11844.1362
$ws.Range("N129").Value = -21844.1362

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 48749.5
$ws.Range("I64").Value = 44999.332
$ws.Range("J64").Value = 60000
$ws.Range("K64").Value = 44999.332
$ws.Range("L64").Value = 60000
$ws.Range("M64").Value = -44751.332
$ws.Range("N64").Value = -60496
$ws.Range("H67").Value = 48749.5
$ws.Range("I67").Value = 44999.332
$ws.Range("J67").Value = 60000
$ws.Range("K67").Value = 44999.332
$ws.Range("L67").Value = 60000
$ws.Range("M67").Value = -44141.332
$ws.Range("N67").Value = -61716
$ws.Range("H70").Value = 4462
$ws.Range("I70").Value = 4354.4
$ws.Range("K70").Value = 4354.4
$ws.Range("M70").Value = -4084.4
$ws.Range("H73").Value = 4462
$ws.Range("I73").Value = 4354.4
$ws.Range("K73").Value = 4354.4
$ws.Range("M73").Value = -3418.4
$ws.Range("H102").Value = 1792.1482
$ws.Range("I102").Value = 1334.0731
$ws.Range("J102").Value = 3236.8462
$ws.Range("K102").Value = 1334.0731
$ws.Range("L102").Value = 3236.8462
$ws.Range("M102").Value = 287.9268999999999
$ws.Range("N102").Value = -6480.8462

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H56").Value = 19147
$ws.Range("I56").Value = 19147
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 19147
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -18456
$ws.Range("N56").ClearContents()
$ws.Range("H132").Value = 117650620
$ws.Range("J132").Value = 250003860
$ws.Range("L132").Value = 750011580
$ws.Range("N132").Value = -750016640

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 24747.5
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H62").Value = 7730.5625
$ws.Range("I62").Value = 8333.333000000001
$ws.Range("J62").Value = 7591.4614
$ws.Range("K62").Value = 8333.333000000001
$ws.Range("L62").Value = 7591.4614
$ws.Range("M62").Value = -7709.333000000001
$ws.Range("N62").Value = -8839.4614
$ws.Range("H65").Value = 7730.5625
$ws.Range("I65").Value = 8333.333000000001
$ws.Range("J65").Value = 7591.4614
$ws.Range("K65").Value = 41666.665
$ws.Range("L65").Value = 37957.307
$ws.Range("M65").Value = -38546.665
$ws.Range("N65").Value = -44197.307
$ws.Range("H96").Value = 7368
$ws.Range("I96").Value = 5191.1665
$ws.Range("J96").Value = 9980.200000000001
$ws.Range("K96").Value = 5191.1665
$ws.Range("L96").Value = 9980.200000000001
$ws.Range("M96").Value = -3818.1665
$ws.Range("N96").Value = -12726.2
$ws.Range("H114").Value = 5000
$ws.Range("J114").Value = 5000
$ws.Range("L114").Value = 5000
$ws.Range("N114").Value = -13678
$ws.Range("H122").Value = 77001864
$ws.Range("I122").Value = 91001384
$ws.Range("K122").Value = 273004152
$ws.Range("M122").Value = -273001702
$ws.Range("H126").Value = 3374.6128
$ws.Range("I126").Value = 3374.6128
$ws.Range("K126").Value = 10123.8384
$ws.Range("M126").Value = -7653.838400000001
$ws.Range("H136").Value = 1413.125
$ws.Range("I136").Value = 1176.0333
$ws.Range("K136").Value = 3528.0999
$ws.Range("M136").Value = -978.0999000000002

